# Adds the new claim/ticket record (row 63) to the "AYKO" sheet, mirroring
# the existing rows: columns A-L are text values, M/N are numeric coordinates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AYKO")

$row = 63

# Force text formatting on A:L first so values that look numeric/date-like
# (e.g. "6221", "9", "1", "6/26/2025") are stored as text, matching the
# workbook's existing inline-string / text cells instead of being
# auto-converted to numbers or date serials.
$textRange = $ws.Range("A$row" + ":L$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = "6221"
$ws.Cells.Item($row, 2).Value  = "6/26/2025"
$ws.Cells.Item($row, 3).Value  = "FERNANDEZ 1549"
$ws.Cells.Item($row, 4).Value  = "9"
$ws.Cells.Item($row, 5).Value  = "807789699"
$ws.Cells.Item($row, 6).Value  = "AYKO"
$ws.Cells.Item($row, 7).Value  = "Pendiente"
$ws.Cells.Item($row, 8).Value  = "Poste inclinado"
$ws.Cells.Item($row, 9).Value  = "1"
$ws.Cells.Item($row, 10).Value = "Aplomo"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Poste"

# M/N are the numeric X/Y coordinates, same as all other rows.
$ws.Cells.Item($row, 13).Value = -58.471717
$ws.Cells.Item($row, 14).Value = -34.649961

# Drop the explicit "@" text style we applied above so the new row's cells
# end up unstyled, same as the other data rows in the sheet.
$textRange.ClearFormats()
